# Actualizacion Datos Personales 4 nov
$wb = $excel.ActiveWorkbook

# --- Rescatables sheet: remove two student rows ---
$rescatables = $wb.Worksheets.Item("Rescatables")
# Delete from bottom up so row numbers of the earlier row stay valid
$rescatables.Rows.Item(18).Delete()
$rescatables.Rows.Item(15).Delete()

# --- Estadisticos 1P sheet: update stats for 1BM (row3) and 3ARHM (row8) ---
$p1 = $wb.Worksheets.Item("Estadisticos 1P")
$p1.Range("D3").Value = 13
$p1.Range("F3").Value = 18
$p1.Range("G3").Value = 58.06
$p1.Range("H3").Value = 7.2

$p1.Range("D8").Value = 5
$p1.Range("F8").Value = 31
$p1.Range("G8").Value = 86.11
$p1.Range("H8").Value = 8.199999999999999

# --- Estadisticos 2P sheet: update Reprobados for 1BM (row3) and 3ARHM (row8) ---
$p2 = $wb.Worksheets.Item("Estadisticos 2P")
$p2.Range("E3").Value = 18
$p2.Range("E8").Value = 31

# --- Estadisticos Final sheet: same stat updates as Estadisticos 1P ---
$pf = $wb.Worksheets.Item("Estadisticos Final")
$pf.Range("D3").Value = 13
$pf.Range("F3").Value = 18
$pf.Range("G3").Value = 58.06
$pf.Range("H3").Value = 7.2

$pf.Range("D8").Value = 5
$pf.Range("F8").Value = 31
$pf.Range("G8").Value = 86.11
$pf.Range("H8").Value = 8.199999999999999
